# Refresh the crypto price/volume figures (columns D and E) for rows 2-51.
# Column D values are plain numeric-looking strings in the source data (e.g. "205.47")
# but must stay TEXT cells (original file stores them as inline strings), so we assign
# them with a leading apostrophe (Excel text-entry prefix) and then reset the range style
# back to "Normal" so no stray number-format / quote-prefix styling is left behind.
# Column E values already contain padding spaces and a "%" sign, so Excel stores them as
# text automatically - no special handling required there.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.662.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "'1.530.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.71%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'205.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("E6").Value = "  -1.08%  "
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "'21.32"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.63%  "
$ws.Range("E9").Value = "  -1.31%  "
$ws.Range("E10").Value = "  -0.69%  "
$ws.Range("D11").Value = "'0.0852"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.37%  "
$ws.Range("D12").Value = "'1.747.02"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.88%  "
$ws.Range("D13").Value = "'1.529.27"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.89%  "
$ws.Range("D14").Value = "'3.66"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.96%  "
$ws.Range("E15").Value = "  -1.16%  "
$ws.Range("D16").Value = "'61.37"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").Value = "'26.660.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").Value = "'212.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.73%  "
$ws.Range("E19").Value = "  +1.34%  "
$ws.Range("D20").Value = "'7.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.75%  "
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("E22").Value = "  -2.08%  "
$ws.Range("D23").Value = "'9.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.15%  "
$ws.Range("D24").Value = "'1.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.51%  "
$ws.Range("D25").Value = "'151.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.83%  "
$ws.Range("E26").Value = "  -3.54%  "
$ws.Range("D27").Value = "'14.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E29").Value = "  -0.75%  "
$ws.Range("E30").Value = "  -0.64%  "
$ws.Range("E31").Value = "  -1.95%  "
$ws.Range("E32").Value = "  +2.76%  "
$ws.Range("D33").Value = "'1.351.28"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.21%  "
$ws.Range("D34").Value = "'2.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.26%  "
$ws.Range("E35").Value = "  -3.85%  "
$ws.Range("D36").Value = "'0.947"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.17%  "
$ws.Range("E37").Value = "  -0.78%  "
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("E39").Value = "  +0.50%  "
$ws.Range("E40").Value = "  -1.53%  "
$ws.Range("D41").Value = "'5.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.20%  "
$ws.Range("D42").Value = "'0.993"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D44").Value = "'62.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.18%  "
$ws.Range("E45").Value = "  -2.22%  "
$ws.Range("D46").Value = "'2.24"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.54%  "
$ws.Range("D47").Value = "'1.662.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.90%  "
$ws.Range("D48").Value = "'85.50"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").Value = "'0.0507"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.91%  "
$ws.Range("D50").Value = "'0.0₇0966"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.40%  "
$ws.Range("D51").Value = "'0.0942"
$ws.Range("D51").Style = "Normal"
